$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right after the header (pushing the existing data rows down)
$ws.Rows("2:4").Insert()

# Populate the newly inserted rows with the new accelerometer readings
$ws.Range("A2").Value = -0.4886150360107422
$ws.Range("B2").Value = 1.498652458190918
$ws.Range("C2").Value = -0.1321379840373993

$ws.Range("A3").Value = -0.7675657272338867
$ws.Range("B3").Value = 1.561143398284912
$ws.Range("C3").Value = -0.3004561066627502

$ws.Range("A4").Value = -0.6989822387695312
$ws.Range("B4").Value = 1.441655874252319
$ws.Range("C4").Value = -0.3177179098129272

# The Insert() above copies the header's formatting down onto the new rows;
# clear it so the new data rows stay unstyled like the rest of the data rows.
$ws.Range("A2:C4").ClearFormats()

# After the insert, the former last four data rows (old rows 19-22) now sit at
# rows 22-25. Drop them so the sheet ends at row 21.
$ws.Rows("22:25").Delete()
